$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156
$ws.Cells.Item(156, 2).Value = 7211640
$ws.Cells.Item(156, 6).Value = "UTC Cajamarca"
$ws.Cells.Item(156, 7).Value = "Sport Boys"
$ws.Cells.Item(156, 8).Value = 1
$ws.Cells.Item(156, 9).Value = 1
$ws.Cells.Item(156, 10).Value = "D"
$ws.Cells.Item(156, 11).Value = 1.615
$ws.Cells.Item(156, 12).Value = 3.75
$ws.Cells.Item(156, 13).Value = 5
$ws.Cells.Item(156, 14).Value = 1.5
$ws.Cells.Item(156, 15).Value = 4.2
$ws.Cells.Item(156, 16).Value = 6.5
$ws.Cells.Item(156, 17).Value = -1
$ws.Cells.Item(156, 18).Value = 1.8
$ws.Cells.Item(156, 19).Value = 2.05
$ws.Cells.Item(156, 20).Value = 2.5
$ws.Cells.Item(156, 21).Value = 1.875
$ws.Cells.Item(156, 22).Value = 1.975
$ws.Cells.Item(156, 23).Value = -1
$ws.Cells.Item(156, 24).Value = 3.2
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = 1.05
$ws.Cells.Item(156, 28).Value = -1
$ws.Cells.Item(156, 29).Value = 0.9750000000000001

# Row 157
$ws.Cells.Item(157, 2).Value = 7211641
$ws.Cells.Item(157, 6).Value = "Sport Huancayo"
$ws.Cells.Item(157, 7).Value = "Deportivo Municipal"
$ws.Cells.Item(157, 8).Value = 2
$ws.Cells.Item(157, 9).Value = 0
$ws.Cells.Item(157, 10).Value = "H"
$ws.Cells.Item(157, 11).Value = 1.125
$ws.Cells.Item(157, 12).Value = 7
$ws.Cells.Item(157, 13).Value = 17
$ws.Cells.Item(157, 14).Value = 1.166
$ws.Cells.Item(157, 15).Value = 6.5
$ws.Cells.Item(157, 16).Value = 12
$ws.Cells.Item(157, 17).Value = -2
$ws.Cells.Item(157, 18).Value = 1.775
$ws.Cells.Item(157, 19).Value = 2.025
$ws.Cells.Item(157, 20).Value = 3.5
$ws.Cells.Item(157, 21).Value = 1.9
$ws.Cells.Item(157, 22).Value = 1.9
$ws.Cells.Item(157, 23).Value = 0.1659999999999999
$ws.Cells.Item(157, 24).Value = -1
$ws.Cells.Item(157, 25).Value = -1
$ws.Cells.Item(157, 26).Value = 0
$ws.Cells.Item(157, 27).Value = -0
$ws.Cells.Item(157, 28).Value = -1
$ws.Cells.Item(157, 29).Value = 0.8999999999999999

# Row 175
$ws.Cells.Item(175, 2).Value = 7302796
$ws.Cells.Item(175, 6).Value = "Sport Huancayo"
$ws.Cells.Item(175, 7).Value = "Sport Boys"
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(175, 9).Value = 0
$ws.Cells.Item(175, 10).Value = "H"
$ws.Cells.Item(175, 11).Value = 1.727
$ws.Cells.Item(175, 12).Value = 3.75
$ws.Cells.Item(175, 13).Value = 4.333
$ws.Cells.Item(175, 14).Value = 1.25
$ws.Cells.Item(175, 15).Value = 5.25
$ws.Cells.Item(175, 16).Value = 10
$ws.Cells.Item(175, 17).Value = -1.75
$ws.Cells.Item(175, 18).Value = 1.925
$ws.Cells.Item(175, 19).Value = 1.875
$ws.Cells.Item(175, 20).Value = 3
$ws.Cells.Item(175, 21).Value = 1.875
$ws.Cells.Item(175, 22).Value = 1.925
$ws.Cells.Item(175, 23).Value = 0.25
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = -1
$ws.Cells.Item(175, 26).Value = -1
$ws.Cells.Item(175, 27).Value = 0.875
$ws.Cells.Item(175, 28).Value = -1
$ws.Cells.Item(175, 29).Value = 0.925

# Row 176
$ws.Cells.Item(176, 2).Value = 7302795
$ws.Cells.Item(176, 6).Value = "Unin Comercio"
$ws.Cells.Item(176, 7).Value = "Deportivo Garcilaso"
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = "A"
$ws.Cells.Item(176, 11).Value = 2.25
$ws.Cells.Item(176, 12).Value = 3.3
$ws.Cells.Item(176, 13).Value = 2.7
$ws.Cells.Item(176, 14).Value = 1.75
$ws.Cells.Item(176, 15).Value = 3.6
$ws.Cells.Item(176, 16).Value = 4
$ws.Cells.Item(176, 17).Value = -0.5
$ws.Cells.Item(176, 18).Value = 1.8
$ws.Cells.Item(176, 19).Value = 2
$ws.Cells.Item(176, 20).Value = 2.75
$ws.Cells.Item(176, 21).Value = 1.825
$ws.Cells.Item(176, 22).Value = 1.975
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 3
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 1
$ws.Cells.Item(176, 28).Value = 0.4125
$ws.Cells.Item(176, 29).Value = -0.5

# Row 177
$ws.Cells.Item(177, 2).Value = 7302200
$ws.Cells.Item(177, 6).Value = "Carlos Manucci"
$ws.Cells.Item(177, 7).Value = "Deportivo Binacional"
$ws.Cells.Item(177, 8).Value = 3
$ws.Cells.Item(177, 9).Value = 2
$ws.Cells.Item(177, 10).Value = "H"
$ws.Cells.Item(177, 11).Value = 2
$ws.Cells.Item(177, 12).Value = 3.2
$ws.Cells.Item(177, 13).Value = 3.75
$ws.Cells.Item(177, 14).Value = 1.75
$ws.Cells.Item(177, 15).Value = 3.4
$ws.Cells.Item(177, 16).Value = 4.333
$ws.Cells.Item(177, 17).Value = -0.5
$ws.Cells.Item(177, 18).Value = 1.85
$ws.Cells.Item(177, 19).Value = 1.95
$ws.Cells.Item(177, 20).Value = 2.5
$ws.Cells.Item(177, 21).Value = 1.85
$ws.Cells.Item(177, 22).Value = 1.95
$ws.Cells.Item(177, 23).Value = 0.75
$ws.Cells.Item(177, 24).Value = -1
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = 0.8500000000000001
$ws.Cells.Item(177, 27).Value = -1
$ws.Cells.Item(177, 28).Value = 0.8500000000000001
$ws.Cells.Item(177, 29).Value = -1

# Row 181
$ws.Cells.Item(181, 2).Value = 7384622
$ws.Cells.Item(181, 6).Value = "Deportivo Municipal"
$ws.Cells.Item(181, 7).Value = "Academia Deportiva Cantolao"
$ws.Cells.Item(181, 8).Value = 1
$ws.Cells.Item(181, 9).Value = 2
$ws.Cells.Item(181, 10).Value = "A"
$ws.Cells.Item(181, 11).Value = 1.444
$ws.Cells.Item(181, 12).Value = 4.333
$ws.Cells.Item(181, 13).Value = 7
$ws.Cells.Item(181, 14).Value = 1.5
$ws.Cells.Item(181, 15).Value = 3.75
$ws.Cells.Item(181, 16).Value = 6
$ws.Cells.Item(181, 17).Value = -1
$ws.Cells.Item(181, 18).Value = 1.825
$ws.Cells.Item(181, 19).Value = 2.025
$ws.Cells.Item(181, 20).Value = 2.75
$ws.Cells.Item(181, 21).Value = 1.875
$ws.Cells.Item(181, 22).Value = 1.975
$ws.Cells.Item(181, 23).Value = -1
$ws.Cells.Item(181, 24).Value = -1
$ws.Cells.Item(181, 25).Value = 5
$ws.Cells.Item(181, 26).Value = -1
$ws.Cells.Item(181, 27).Value = 1.025
$ws.Cells.Item(181, 28).Value = 0.4375
$ws.Cells.Item(181, 29).Value = -0.5

# Row 182
$ws.Cells.Item(182, 2).Value = 7384623
$ws.Cells.Item(182, 6).Value = "Sport Boys"
$ws.Cells.Item(182, 7).Value = "Cienciano"
$ws.Cells.Item(182, 8).Value = 2
$ws.Cells.Item(182, 9).Value = 1
$ws.Cells.Item(182, 10).Value = "H"
$ws.Cells.Item(182, 11).Value = 2
$ws.Cells.Item(182, 12).Value = 3.4
$ws.Cells.Item(182, 13).Value = 3.5
$ws.Cells.Item(182, 14).Value = 1.833
$ws.Cells.Item(182, 15).Value = 4
$ws.Cells.Item(182, 16).Value = 3.2
$ws.Cells.Item(182, 17).Value = -0.5
$ws.Cells.Item(182, 18).Value = 1.925
$ws.Cells.Item(182, 19).Value = 1.875
$ws.Cells.Item(182, 20).Value = 3
$ws.Cells.Item(182, 21).Value = 1.925
$ws.Cells.Item(182, 22).Value = 1.875
$ws.Cells.Item(182, 23).Value = 0.833
$ws.Cells.Item(182, 24).Value = -1
$ws.Cells.Item(182, 25).Value = -1
$ws.Cells.Item(182, 26).Value = 0.925
$ws.Cells.Item(182, 27).Value = -1
$ws.Cells.Item(182, 28).Value = 0
$ws.Cells.Item(182, 29).Value = -0

# Row 183
$ws.Cells.Item(183, 2).Value = 7384628
$ws.Cells.Item(183, 6).Value = "Deportivo Binacional"
$ws.Cells.Item(183, 7).Value = "FBC Melgar"
$ws.Cells.Item(183, 8).Value = 1
$ws.Cells.Item(183, 9).Value = 2
$ws.Cells.Item(183, 10).Value = "A"
$ws.Cells.Item(183, 11).Value = 2.75
$ws.Cells.Item(183, 12).Value = 3.3
$ws.Cells.Item(183, 13).Value = 2.375
$ws.Cells.Item(183, 14).Value = 3.3
$ws.Cells.Item(183, 15).Value = 3.6
$ws.Cells.Item(183, 16).Value = 2
$ws.Cells.Item(183, 17).Value = 0.5
$ws.Cells.Item(183, 18).Value = 1.8
$ws.Cells.Item(183, 19).Value = 2
$ws.Cells.Item(183, 20).Value = 2.75
$ws.Cells.Item(183, 21).Value = 1.975
$ws.Cells.Item(183, 22).Value = 1.875
$ws.Cells.Item(183, 23).Value = -1
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 25).Value = 1
$ws.Cells.Item(183, 26).Value = -1
$ws.Cells.Item(183, 27).Value = 1
$ws.Cells.Item(183, 28).Value = 0.4875
$ws.Cells.Item(183, 29).Value = -0.5

# Row 184
$ws.Cells.Item(184, 2).Value = 7384630
$ws.Cells.Item(184, 6).Value = "Atletico Grau"
$ws.Cells.Item(184, 7).Value = "Unin Comercio"
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = 1
$ws.Cells.Item(184, 10).Value = "A"
$ws.Cells.Item(184, 11).Value = 2.8
$ws.Cells.Item(184, 12).Value = 3.4
$ws.Cells.Item(184, 13).Value = 2.15
$ws.Cells.Item(184, 14).Value = 1.75
$ws.Cells.Item(184, 15).Value = 3.6
$ws.Cells.Item(184, 16).Value = 3.8
$ws.Cells.Item(184, 17).Value = -0.75
$ws.Cells.Item(184, 18).Value = 2
$ws.Cells.Item(184, 19).Value = 1.8
$ws.Cells.Item(184, 20).Value = 3
$ws.Cells.Item(184, 21).Value = 1.85
$ws.Cells.Item(184, 22).Value = 1.95
$ws.Cells.Item(184, 23).Value = -1
$ws.Cells.Item(184, 24).Value = -1
$ws.Cells.Item(184, 25).Value = 2.8
$ws.Cells.Item(184, 26).Value = -1
$ws.Cells.Item(184, 27).Value = 0.8
$ws.Cells.Item(184, 28).Value = -1
$ws.Cells.Item(184, 29).Value = 0.95

# Row 185
$ws.Cells.Item(185, 2).Value = 7384625
$ws.Cells.Item(185, 6).Value = "AD Tarma"
$ws.Cells.Item(185, 7).Value = "Carlos Manucci"
$ws.Cells.Item(185, 8).Value = 0
$ws.Cells.Item(185, 9).Value = 0
$ws.Cells.Item(185, 10).Value = "D"
$ws.Cells.Item(185, 11).Value = 1.5
$ws.Cells.Item(185, 12).Value = 3.75
$ws.Cells.Item(185, 13).Value = 7
$ws.Cells.Item(185, 14).Value = 1.363
$ws.Cells.Item(185, 15).Value = 4.333
$ws.Cells.Item(185, 16).Value = 9.5
$ws.Cells.Item(185, 17).Value = -1.25
$ws.Cells.Item(185, 18).Value = 1.875
$ws.Cells.Item(185, 19).Value = 1.925
$ws.Cells.Item(185, 20).Value = 2.5
$ws.Cells.Item(185, 21).Value = 1.8
$ws.Cells.Item(185, 22).Value = 2
$ws.Cells.Item(185, 23).Value = -1
$ws.Cells.Item(185, 24).Value = 3.333
$ws.Cells.Item(185, 25).Value = -1
$ws.Cells.Item(185, 26).Value = -1
$ws.Cells.Item(185, 27).Value = 0.925
$ws.Cells.Item(185, 28).Value = -1
$ws.Cells.Item(185, 29).Value = 1

# Row 186
$ws.Cells.Item(186, 2).Value = 7384626
$ws.Cells.Item(186, 6).Value = "Sporting Cristal"
$ws.Cells.Item(186, 7).Value = "Alianza Atletico"
$ws.Cells.Item(186, 8).Value = 3
$ws.Cells.Item(186, 9).Value = 0
$ws.Cells.Item(186, 10).Value = "H"
$ws.Cells.Item(186, 11).Value = 1.3
$ws.Cells.Item(186, 12).Value = 5
$ws.Cells.Item(186, 13).Value = 9
$ws.Cells.Item(186, 14).Value = 1.166
$ws.Cells.Item(186, 15).Value = 6.5
$ws.Cells.Item(186, 16).Value = 13
$ws.Cells.Item(186, 17).Value = -2
$ws.Cells.Item(186, 18).Value = 1.85
$ws.Cells.Item(186, 19).Value = 1.95
$ws.Cells.Item(186, 20).Value = 3.25
$ws.Cells.Item(186, 21).Value = 2
$ws.Cells.Item(186, 22).Value = 1.8
$ws.Cells.Item(186, 23).Value = 0.1659999999999999
$ws.Cells.Item(186, 24).Value = -1
$ws.Cells.Item(186, 25).Value = -1
$ws.Cells.Item(186, 26).Value = 0.8500000000000001
$ws.Cells.Item(186, 27).Value = -1
$ws.Cells.Item(186, 28).Value = -0.5
$ws.Cells.Item(186, 29).Value = 0.4

# Row 187
$ws.Cells.Item(187, 2).Value = 7384629
$ws.Cells.Item(187, 6).Value = "Deportivo Garcilaso"
$ws.Cells.Item(187, 7).Value = "Alianza Lima"
$ws.Cells.Item(187, 8).Value = 0
$ws.Cells.Item(187, 9).Value = 1
$ws.Cells.Item(187, 10).Value = "A"
$ws.Cells.Item(187, 11).Value = 2.625
$ws.Cells.Item(187, 12).Value = 3.3
$ws.Cells.Item(187, 13).Value = 2.5
$ws.Cells.Item(187, 14).Value = 2.7
$ws.Cells.Item(187, 15).Value = 3.4
$ws.Cells.Item(187, 16).Value = 2.375
$ws.Cells.Item(187, 17).Value = 0
$ws.Cells.Item(187, 18).Value = 2.025
$ws.Cells.Item(187, 19).Value = 1.775
$ws.Cells.Item(187, 20).Value = 2.25
$ws.Cells.Item(187, 21).Value = 1.825
$ws.Cells.Item(187, 22).Value = 1.975
$ws.Cells.Item(187, 23).Value = -1
$ws.Cells.Item(187, 24).Value = -1
$ws.Cells.Item(187, 25).Value = 1.375
$ws.Cells.Item(187, 26).Value = -1
$ws.Cells.Item(187, 27).Value = 0.7749999999999999
$ws.Cells.Item(187, 28).Value = -1
$ws.Cells.Item(187, 29).Value = 0.9750000000000001

# Row 188
$ws.Cells.Item(188, 2).Value = 7384627
$ws.Cells.Item(188, 6).Value = "Universitario de Deportes"
$ws.Cells.Item(188, 7).Value = "Sport Huancayo"
$ws.Cells.Item(188, 8).Value = 2
$ws.Cells.Item(188, 9).Value = 0
$ws.Cells.Item(188, 10).Value = "H"
$ws.Cells.Item(188, 11).Value = 1.25
$ws.Cells.Item(188, 12).Value = 5
$ws.Cells.Item(188, 13).Value = 12
$ws.Cells.Item(188, 14).Value = 1.181
$ws.Cells.Item(188, 15).Value = 6
$ws.Cells.Item(188, 16).Value = 13
$ws.Cells.Item(188, 17).Value = -1.75
$ws.Cells.Item(188, 18).Value = 1.8
$ws.Cells.Item(188, 19).Value = 2
$ws.Cells.Item(188, 20).Value = 2.75
$ws.Cells.Item(188, 21).Value = 1.85
$ws.Cells.Item(188, 22).Value = 1.95
$ws.Cells.Item(188, 23).Value = 0.181
$ws.Cells.Item(188, 24).Value = -1
$ws.Cells.Item(188, 25).Value = -1
$ws.Cells.Item(188, 26).Value = 0.4
$ws.Cells.Item(188, 27).Value = -0.5
$ws.Cells.Item(188, 28).Value = -1
$ws.Cells.Item(188, 29).Value = 0.95
